# Remove the "Корр. счет: @<CORR_ACCOUNT>@" paragraph entirely (wrong/duplicate
# tag removed per commit message "wrong tag and remove cor. account").
# The preceding "БИК: @<BIK>@" paragraph is left untouched.

$d = $word.ActiveDocument

$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*CORR_ACCOUNT*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
